$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the model-name labels in column A (rows 2-11) with more interpretable
# names (full model descriptions instead of abbreviations). Order of
# assignment matches the order the original author entered/edited them in.
$ws.Range("A4").Value  = "year + site"
$ws.Range("A6").Value  = "year + site + year*site"
$ws.Range("A11").Value = "year"
$ws.Range("A3").Value  = "rugosity + site"
$ws.Range("A5").Value  = "rugosity + year + site"
$ws.Range("A7").Value  = "rugosity + site + site*rugosity"
$ws.Range("A8").Value  = "rugosity + year + year*rugosity"
$ws.Range("A9").Value  = "rugosity + year"
$ws.Range("A10").Value = "rugosity"

# Widen column A to fit the new, longer model names (re-running "best fit"
# now that the labels are longer).
$ws.Columns.Item(1).ColumnWidth = 27.83

# Clear the selection so the sheet view no longer carries the previous range
# selection.
$ws.Range("A1").Select()
